$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.704.46"
$ws.Range("E2").Value = "  -1.82%  "
$ws.Range("D3").Value = "1.801.93"
$ws.Range("E3").Value = "  -1.48%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.43"
$ws.Range("E5").Value = "  -2.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5950"
$ws.Range("E6").Value = "  -2.59%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2776"
$ws.Range("E8").Value = "  -1.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06825"
$ws.Range("E9").Value = "  -3.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.31"
$ws.Range("E10").Value = "  -1.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07548"
$ws.Range("E11").Value = "  -1.68%  "
$ws.Range("D12").Value = "1.802.96"
$ws.Range("E12").Value = "  -1.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.761"
$ws.Range("E13").Value = "  -1.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6229"
$ws.Range("D15").Value = "2.047.45"
$ws.Range("E15").Value = "  -1.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009298"
$ws.Range("E16").Value = "  -7.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "75.38"
$ws.Range("E17").Value = "  -4.33%  "
$ws.Range("D18").Value = "28.647.34"
$ws.Range("E18").Value = "  -2.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.477"
$ws.Range("E19").Value = "  -6.52%  "
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "209.66"
$ws.Range("E21").Value = "  -7.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.43"
$ws.Range("E22").Value = "  -3.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.847"
$ws.Range("E23").Value = "  -2.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.003"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.19"
$ws.Range("E25").Value = "  -0.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.840"
$ws.Range("E26").Value = "  -2.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1272"
$ws.Range("E27").Value = "  -2.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.39"
$ws.Range("E28").Value = "  -1.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.427"
$ws.Range("E29").Value = "  -3.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06162"
$ws.Range("E30").Value = "  -3.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.776"
$ws.Range("E32").Value = "  -1.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.739"
$ws.Range("E33").Value = "  -1.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.719"
$ws.Range("E34").Value = "  -1.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.063"
$ws.Range("E35").Value = "  -5.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6382"
$ws.Range("E36").Value = "  -1.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.491"
$ws.Range("E37").Value = "  -2.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.716"
$ws.Range("E38").Value = "  -0.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.451"
$ws.Range("E39").Value = "  -1.84%  "
$ws.Range("E40").Value = "  -1.76%  "
$ws.Range("D41").Value = "1.131.18"
$ws.Range("E41").Value = "  -6.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8719"
$ws.Range("E42").Value = "  -4.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.004"
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.67"
$ws.Range("E44").Value = "  -0.41%  "
$ws.Range("D45").Value = "1.964.15"
$ws.Range("E45").Value = "  -0.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.48"
$ws.Range("E46").Value = "  -3.76%  "
$ws.Range("E47").Value = "  -3.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.598"
$ws.Range("E48").Value = "  -1.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05463"
$ws.Range("E49").Value = "  -1.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.316"
$ws.Range("E50").Value = "  -3.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4489"
$ws.Range("E51").Value = "  -1.75%  "
